$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto prices / volume percentages, and the two
# coin re-rankings (rows 32/33 and 42/43 swap content) from the feed refresh.

$ws.Range("D2").Value = "'71.973.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.41%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'4.004.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.20%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'532.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.92%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'152.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.38%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.695"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +10.73%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.745"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.84%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  -3.71%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  -6.07%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'47.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +6.38%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'10.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.23%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'4.644.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.00%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'4.010.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.71%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'13.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.91%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'20.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -4.43%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  -0.90%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'1.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.06%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'71.882.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.35%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'425.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.64%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'97.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +3.95%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  -2.41%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  +2.29%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'14.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.65%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'11.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -9.74%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'10.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.80%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +1.32%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'36.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.54%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'3.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +24.25%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'13.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.22%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").Value = "'NEARProtocol"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'7.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.76%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").Value = "'Hedera"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'0.129"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.72%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'670.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.88%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'65.78"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.10%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'42.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.55%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.426"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -4.99%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +0.28%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.0₃0824"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -9.63%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'3.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -6.32%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.24%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "'FirstDigitalUSD"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.00%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = "'WEMIXToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'3.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +4.34%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  -1.95%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'  +2.24%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'9.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +3.47%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  -9.76%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  -5.01%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'2.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -7.83%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.000272"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -4.31%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'144.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.82%  "
$ws.Range("E51").Style = "Normal"
